$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 9-11 (column D and E change; F stays the same)
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 3

$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 6

$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 7

# Row 13: D/E change, F stays 0
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 7

# New rows 14-16
# Copy A13's cell format (bold/border/center) down to the new A-column cells
$ws.Range("A13").Copy()
$ws.Range("A14:A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 3

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 17
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 6

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 21
$ws.Range("E16").Value = 12
$ws.Range("F16").Value = 7
